$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_summary")

$ws.Range("B2").Value = -14400.73170010763
$ws.Range("C2").Value = 30035.40918737522
$ws.Range("D2").Value = 29165.46340021527

$ws.Range("B3").Value = -14244.73077293404
$ws.Range("C3").Value = 30136.98256938385
$ws.Range("D3").Value = 28975.46154586809

$ws.Range("B4").Value = -14106.66270581043
$ws.Range("C4").Value = 30274.42167149243
$ws.Range("D4").Value = 28821.32541162087

$ws.Range("B5").Value = -13994.28327716091
$ws.Range("C5").Value = 30463.2380505492
$ws.Range("D5").Value = 28718.56655432183
